$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (45171 -> 45172) for every data row (rows 2 through 496) as part of
# an automatic update.
$ws.Range("C2:C496").Value = 45172
